# Auto-applied scheduled market-data refresh for Leve profit sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# with freshly fetched values for the rows whose item prices changed.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 243.63637
$ws.Range("I33").Value = 258.8
$ws.Range("K33").Value = 258.8
$ws.Range("M33").Value = -29.80000000000001
$ws.Range("H41").Value = 271.63635
$ws.Range("I41").Value = 271.63635
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 271.63635
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = 168.36365
$ws.Range("H53").Value = 278.2353
$ws.Range("I53").Value = 229.08333
$ws.Range("J53").Value = 396.2
$ws.Range("K53").Value = 229.08333
$ws.Range("L53").Value = 396.2
$ws.Range("M53").Value = 407.91667
$ws.Range("N53").Value = -1670.2
$ws.Range("H86").Value = 2874.25
$ws.Range("I86").Value = 2832.3333
$ws.Range("K86").Value = 2832.3333
$ws.Range("M86").Value = -1709.3333
$ws.Range("H89").Value = 2874.25
$ws.Range("I89").Value = 2832.3333
$ws.Range("K89").Value = 14161.6665
$ws.Range("M89").Value = -8545.666499999999
$ws.Range("H111").Value = 1400
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H137").Value = 5475.923
$ws.Range("I137").Value = 2993.5
$ws.Range("K137").Value = 8980.5
$ws.Range("M137").Value = -6430.5
$ws.Range("H138").Value = 4301.0938
$ws.Range("J138").Value = 4519.148
$ws.Range("L138").Value = 13557.444
$ws.Range("N138").Value = -23837.444

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2813.5588
$ws.Range("I32").Value = 1895.6875
$ws.Range("K32").Value = 1895.6875
$ws.Range("M32").Value = -1608.6875
$ws.Range("H63").Value = 3768.75
$ws.Range("I63").Value = 3768.75
$ws.Range("K63").Value = 3768.75
$ws.Range("M63").Value = -3082.75
$ws.Range("H66").Value = 3768.75
$ws.Range("I66").Value = 3768.75
$ws.Range("K66").Value = 18843.75
$ws.Range("M66").Value = -15411.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 26036.666
$ws.Range("I82").Value = 3913.5
$ws.Range("K82").Value = 3913.5
$ws.Range("M82").Value = -3530.5
$ws.Range("H85").Value = 26036.666
$ws.Range("I85").Value = 3913.5
$ws.Range("K85").Value = 3913.5
$ws.Range("M85").Value = -2587.5
$ws.Range("H94").Value = 1342
$ws.Range("I94").Value = 1110
$ws.Range("J94").Value = 1574
$ws.Range("K94").Value = 1110
$ws.Range("L94").Value = 1574
$ws.Range("M94").Value = -659
$ws.Range("N94").Value = -2476
$ws.Range("H99").Value = 789.5
$ws.Range("I99").Value = 789.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 789.5
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 708.5
$ws.Range("H107").Value = 8233.333000000001
$ws.Range("I107").Value = 8233.333000000001
$ws.Range("K107").Value = 8233.333000000001
$ws.Range("M107").Value = -6313.333000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2122.3333
$ws.Range("I58").Value = 2122.3333
$ws.Range("K58").Value = 2122.3333
$ws.Range("M58").Value = -1919.3333
$ws.Range("H68").Value = 70295
$ws.Range("J68").Value = 70295
$ws.Range("L68").Value = 70295
$ws.Range("N68").Value = -71793
$ws.Range("H71").Value = 70295
$ws.Range("J71").Value = 70295
$ws.Range("L71").Value = 210885
$ws.Range("N71").Value = -218373
$ws.Range("H74").Value = 56189.668
$ws.Range("J74").Value = 56189.668
$ws.Range("L74").Value = 56189.668
$ws.Range("N74").Value = -57937.668
$ws.Range("H77").Value = 56189.668
$ws.Range("J77").Value = 56189.668
$ws.Range("L77").Value = 168569.004
$ws.Range("N77").Value = -177305.004
$ws.Range("H107").Value = 1007.5714
$ws.Range("I107").Value = 588.5
$ws.Range("J107").Value = 1566.3334
$ws.Range("K107").Value = 588.5
$ws.Range("L107").Value = 1566.3334
$ws.Range("M107").Value = 1331.5
$ws.Range("N107").Value = -5406.3334
$ws.Range("H134").Value = 3735.625
$ws.Range("I134").Value = 3735.625
$ws.Range("K134").Value = 11206.875
$ws.Range("M134").Value = -8671.875
$ws.Range("H136").Value = 2122.3333
$ws.Range("I136").Value = 2122.3333
$ws.Range("K136").Value = 6366.999899999999
$ws.Range("M136").Value = -3816.999899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1370.3334
$ws.Range("J129").Value = 1806.5
$ws.Range("L129").Value = 5419.5
$ws.Range("N129").Value = -15419.5
$ws.Range("H131").Value = 734613.6
$ws.Range("J131").Value = 918019.9399999999
$ws.Range("L131").Value = 2754059.82
$ws.Range("N131").Value = -2764139.82
$ws.Range("H132").Value = 3819.6667
$ws.Range("I132").Value = 966
$ws.Range("K132").Value = 8694
$ws.Range("M132").Value = -6164
$ws.Range("H139").Value = 5187
$ws.Range("I139").Value = 3827.25
$ws.Range("K139").Value = 11481.75
$ws.Range("M139").Value = -6341.75
$ws.Range("H140").Value = 2499
$ws.Range("I140").Value = 1248.5
$ws.Range("K140").Value = 3745.5
$ws.Range("M140").Value = 1434.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 4000
$ws.Range("I19").Value = 4000
$ws.Range("K19").Value = 4000
$ws.Range("M19").Value = -3712
$ws.Range("H113").Value = 1993.5
$ws.Range("I113").Value = 1993.5
$ws.Range("K113").Value = 1993.5
$ws.Range("M113").Value = 176.5
$ws.Range("H132").Value = 4820.364
$ws.Range("I132").Value = 4838.8335
$ws.Range("K132").Value = 14516.5005
$ws.Range("M132").Value = -11986.5005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4830.3335
$ws.Range("I40").Value = 4995.5
$ws.Range("J40").Value = 4500
$ws.Range("K40").Value = 4995.5
$ws.Range("L40").Value = 4500
$ws.Range("M40").Value = -4859.5
$ws.Range("N40").Value = -4772
$ws.Range("H46").Value = 3725
$ws.Range("J46").Value = 3737.5
$ws.Range("L46").Value = 3737.5
$ws.Range("N46").Value = -4113.5
$ws.Range("H100").Value = 1674.75
$ws.Range("I100").Value = 1674.75
$ws.Range("K100").Value = 1674.75
$ws.Range("M100").Value = -1133.75
$ws.Range("H133").Value = 135000
$ws.Range("J133").Value = 135000
$ws.Range("L133").Value = 135000
$ws.Range("N133").Value = -140060

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2000
$ws.Range("I96").Value = 2000
$ws.Range("K96").Value = 2000
$ws.Range("M96").Value = -627
$ws.Range("H100").Value = 2298.5
$ws.Range("I100").Value = 2298.5
$ws.Range("K100").Value = 4597
$ws.Range("M100").Value = -4056
$ws.Range("H113").Value = 705
$ws.Range("I113").Value = 698.5714
$ws.Range("K113").Value = 2095.7142
$ws.Range("M113").Value = 74.28579999999965
$ws.Range("H136").Value = 11888.667
$ws.Range("I136").Value = 11888.667
$ws.Range("K136").Value = 35666.001
$ws.Range("M136").Value = -33116.001

